$wb = $excel.ActiveWorkbook

# --- 1. Rename the "research" sheet to "setup" ---
$wsSetup = $wb.Worksheets.Item("research")
$wsSetup.Name = "setup"

$wsResults = $wb.Worksheets.Item("results")

# --- 2. Make "setup" (first sheet) the active sheet / active tab ---
#     (previously the workbook opened with the "results" tab selected;
#      now it should open on "setup", with no explicit topLeftCell and
#      a fresh selection/zoom captured from the editing session)
$wsSetup.Activate()
$excel.ActiveWindow.Zoom = 75
$wsSetup.Range("L14").Select()

# --- 3. Update the experiment note text in M2, now spanning two lines ---
$wsSetup.Range("M2").Value = "Yeni sunucuda ilk uçtan uca deneme yapıldı. `nÖlçüm sadece fold-001 üzerinden gerçekleştirildi."

# Give the merged note cell (M2:M6) a centered vertical alignment and
# wrap the (now multi-line) text so it displays properly.
$noteRange = $wsSetup.Range("M2:M6")
$noteRange.VerticalAlignment = -4108
$wsSetup.Range("M2").WrapText = $true

# Typing the multi-line note auto-expands row 2's height; restore the
# row back to the sheet's default auto-fit height.
$wsSetup.Rows.Item(2).AutoFit()
